$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.369.05"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "1.871.09"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.85"
$ws.Range("E5").Value = "  -1.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4707"
$ws.Range("E7").Value = "  -1.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2881"
$ws.Range("E8").Value = "  -1.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06462"
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.96"
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07778"
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("B12").Value = "Litecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.19"
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.869.90"
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7253"
$ws.Range("E14").Value = "  -2.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.134"
$ws.Range("E15").Value = "  -1.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "280.03"
$ws.Range("E16").Value = "  +1.92%  "
$ws.Range("D17").Value = "30.359.53"
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.03"
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007498"
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("D21").Value = "2.114.30"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.248"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.233"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "163.79"
$ws.Range("E25").Value = "  -1.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.053"
$ws.Range("E26").Value = "  -1.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.71"
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.881"
$ws.Range("E28").Value = "  -1.86%  "
$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09632"
$ws.Range("E29").Value = "  -2.41%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.321"
$ws.Range("E30").Value = "  -1.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.482"
$ws.Range("E31").Value = "  -1.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.233"
$ws.Range("E32").Value = "  -1.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.116"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04815"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.121"
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6902"
$ws.Range("E36").Value = "  -0.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.711"
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01896"
$ws.Range("E38").Value = "  +1.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.811"
$ws.Range("E39").Value = "  +1.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.242"
$ws.Range("E40").Value = "  -0.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.31"
$ws.Range("E41").Value = "  +1.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4227"
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("E43").Value = "  -3.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9993"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8241"
$ws.Range("E45").Value = "  -1.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.83"
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.631"
$ws.Range("E47").Value = "  +2.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.32"
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.958"
$ws.Range("E49").Value = "  -1.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "899.94"
$ws.Range("E50").Value = "  -1.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05723"
$ws.Range("E51").Value = "  +0.63%  "
